# Update "想去人数" (want-to-go count) values in column F across the four
# worksheets, per the regenerated data snapshot (gh-pages output @ 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 (Exhibition) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 2226
$ws1.Range("F3").Value = 252
$ws1.Range("F4").Value = 158
$ws1.Range("F5").Value = 154
$ws1.Range("F6").Value = 291
$ws1.Range("F8").Value = 661
$ws1.Range("F10").Value = 593
$ws1.Range("F11").Value = 361
$ws1.Range("F12").Value = 53
$ws1.Range("F14").Value = 939
$ws1.Range("F17").Value = 79
$ws1.Range("F18").Value = 6
$ws1.Range("F20").Value = 198
$ws1.Range("F21").Value = 76

# --- Sheet 2: 演出 (Performance) ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F3").Value = 29
$ws2.Range("F8").Value = 2387
$ws2.Range("F10").Value = 13
$ws2.Range("F16").Value = 2200

# --- Sheet 3: 本地生活 (Local life) ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F3").Value = 304

# --- Sheet 4: 全部类型 (All types, combined listing) ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F4").Value = 29
$ws4.Range("F5").Value = 2226
$ws4.Range("F6").Value = 304
$ws4.Range("F7").Value = 252
$ws4.Range("F8").Value = 158
$ws4.Range("F9").Value = 154
$ws4.Range("F10").Value = 291
$ws4.Range("F16").Value = 661
$ws4.Range("F18").Value = 593
$ws4.Range("F19").Value = 361
$ws4.Range("F20").Value = 53
$ws4.Range("F22").Value = 939
$ws4.Range("F24").Value = 2388
$ws4.Range("F26").Value = 13
$ws4.Range("F32").Value = 79
$ws4.Range("F33").Value = 6
$ws4.Range("F37").Value = 198
$ws4.Range("F38").Value = 76
$ws4.Range("F39").Value = 2200
